# Update cryptocurrency price/volume data per the latest scrape.
# Cell values are plain numeric-looking text (e.g. "553.75") that must stay
# TEXT, not be auto-converted to numbers -- same as typing them in Excel with
# a leading apostrophe. Build the value and prefix with "'" when needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $looksNumeric = $text -match '^\s*[+-]?((\d+(\.\d*)?)|(\.\d+))\s*$'
    if ($looksNumeric) {
        $cell.Value = "'" + $text
    } else {
        $cell.Value = $text
    }
}

Set-TextValue 2 4 '59.512.64'
Set-TextValue 2 5 '  +1.27%  '

Set-TextValue 3 4 '2.587.32'
Set-TextValue 3 5 '  +0.65%  '

Set-TextValue 4 5 '  -0.06%  '

Set-TextValue 5 4 '553.75'
Set-TextValue 5 5 '  -1.33%  '

Set-TextValue 6 4 '140.37'
Set-TextValue 6 5 '  -1.35%  '

Set-TextValue 7 5 '  -0.17%  '

Set-TextValue 8 5 '  +0.22%  '

Set-TextValue 9 4 '2.607.28'
Set-TextValue 9 5 '  +1.25%  '

Set-TextValue 10 5 '  +0.96%  '

Set-TextValue 11 5 '  +1.78%  '

Set-TextValue 12 5 '  +5.85%  '

Set-TextValue 13 5 '  +5.21%  '

Set-TextValue 14 4 '3.052.45'
Set-TextValue 14 5 '  +0.94%  '

Set-TextValue 15 4 '59.500.02'
Set-TextValue 15 5 '  +1.06%  '

Set-TextValue 16 4 '23.06'
Set-TextValue 16 5 '  +5.61%  '

Set-TextValue 17 5 '  +1.32%  '

Set-TextValue 18 4 '2.600.90'
Set-TextValue 18 5 '  +0.52%  '

Set-TextValue 19 5 '  +1.56%  '

Set-TextValue 20 4 '340.11'
Set-TextValue 20 5 '  +1.63%  '

Set-TextValue 21 4 '10.42'
Set-TextValue 21 5 '  +2.90%  '

Set-TextValue 22 4 '6.59'
Set-TextValue 22 5 '  +7.37%  '

Set-TextValue 23 4 '0.997'
Set-TextValue 23 5 '  -0.22%  '

Set-TextValue 24 4 '0.485'
Set-TextValue 24 5 '  +9.26%  '

Set-TextValue 25 4 '62.79'
Set-TextValue 25 5 '  -1.77%  '

Set-TextValue 26 5 '  -0.51%  '

Set-TextValue 27 5 '  -1.07%  '

Set-TextValue 28 5 '  +4.19%  '

Set-TextValue 29 4 '0.0₃0769'
Set-TextValue 29 5 '  -0.92%  '

Set-TextValue 30 5 '  -0.12%  '

Set-TextValue 31 5 '  +0.77%  '

Set-TextValue 32 4 '6.13'
Set-TextValue 32 5 '  +1.71%  '

Set-TextValue 33 4 '157.51'
Set-TextValue 33 5 '  -1.88%  '

Set-TextValue 34 4 '19.37'
Set-TextValue 34 5 '  +2.74%  '

Set-TextValue 35 4 '4.08'
Set-TextValue 35 5 '  +2.20%  '

Set-TextValue 36 4 '0.919'
Set-TextValue 36 5 '  +5.11%  '

Set-TextValue 37 5 '  +3.32%  '

Set-TextValue 38 5 '  +2.40%  '

Set-TextValue 39 4 '1.48'
Set-TextValue 39 5 '  +0.73%  '

Set-TextValue 40 5 '  -4.03%  '

Set-TextValue 41 5 '  +2.06%  '

Set-TextValue 42 4 '289.79'
Set-TextValue 42 5 '  -1.50%  '

Set-TextValue 43 4 '136.64'
Set-TextValue 43 5 '  +9.13%  '

Set-TextValue 44 4 '0.998'
Set-TextValue 44 5 '  -0.08%  '

Set-TextValue 45 5 '  +0.61%  '

Set-TextValue 46 4 '0.600'
Set-TextValue 46 5 '  +1.20%  '

Set-TextValue 47 2 'Hedera'
Set-TextValue 47 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 47 4 '0.0535'
Set-TextValue 47 5 '  +0.19%  '

Set-TextValue 48 2 'WhiteBITCoin'
Set-TextValue 48 3 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue 48 4 '10.63'
Set-TextValue 48 5 '  +0.14%  '

Set-TextValue 49 5 '  +2.05%  '

Set-TextValue 50 2 'RenderToken'
Set-TextValue 50 3 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 50 4 '4.76'
Set-TextValue 50 5 '  +7.02%  '

Set-TextValue 51 2 'Maker'
Set-TextValue 51 3 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 51 4 '1.970.75'
Set-TextValue 51 5 '  +2.43%  '
